$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.083.62'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '3.755.93'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''602.26'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').Value = '''166.35'
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('D7').Value = '3.754.75'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').Value = '''37.69'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').Value = '4.382.50'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').Value = '3.750.81'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '69.071.38'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = '''7.40'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').Value = '''17.70'
$ws.Range('E19').Value = '  +3.64%  '
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('D21').Value = '''11.22'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('D22').Value = '''491.00'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''84.67'
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '''0.0000148'
$ws.Range('E25').Value = '  -2.01%  '
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('D27').Value = '''12.27'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').Value = '''8.09'
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('D33').Value = '''31.70'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('D34').Value = '3.902.47'
$ws.Range('D35').Value = '3.720.24'
$ws.Range('E35').Value = '  +0.92%  '
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('E37').Value = '  +5.42%  '
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = '''3.09'
$ws.Range('E41').Value = '  +7.46%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '''428.28'
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('D44').Value = '''48.55'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Value = '''40.24'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''142.69'
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('D50').Value = '2.811.85'
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('E51').Value = '  +8.64%  '
